$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "d"
$ws.Range("D2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 48.2625
